$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 35666.668
$ws.Range("J3").Value = 35666.668
$ws.Range("L3").Value = 35666.668
$ws.Range("N3").Value = -35894.668

$ws.Range("H54").Value = 19999
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 19999
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 19999
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -20971

$ws.Range("H63").Value = 40271
$ws.Range("J63").Value = 40271
$ws.Range("L63").Value = 40271
$ws.Range("N63").Value = -41519

$ws.Range("H66").Value = 40271
$ws.Range("J66").Value = 40271
$ws.Range("L66").Value = 120813
$ws.Range("N66").Value = -127053

$ws.Range("H102").Value = 35666.668
$ws.Range("J102").Value = 35666.668
$ws.Range("L102").Value = 35666.668
$ws.Range("N102").Value = -42156.668

$ws.Range("H129").Value = 963.4839
$ws.Range("I129").Value = 321.5
$ws.Range("J129").Value = 1007.7586
$ws.Range("K129").Value = 964.5
$ws.Range("L129").Value = 3023.2758
$ws.Range("M129").Value = 4035.5
$ws.Range("N129").Value = -13023.2758

$ws.Range("H132").Value = 2911.05
$ws.Range("I132").Value = 2736.7856
$ws.Range("J132").Value = 3317.6667
$ws.Range("K132").Value = 8210.356800000001
$ws.Range("L132").Value = 9953.000100000001
$ws.Range("M132").Value = -5680.356800000001
$ws.Range("N132").Value = -15013.0001

$ws.Range("H140").Value = 138865
$ws.Range("J140").Value = 138865
$ws.Range("L140").Value = 138865
$ws.Range("N140").Value = -149225

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 19900
$ws.Range("J13").Value = 19900
$ws.Range("L13").Value = 19900
$ws.Range("N13").Value = -20236

$ws.Range("H54").Value = 5928.222
$ws.Range("I54").Value = 2841.25
$ws.Range("J54").Value = 8397.8
$ws.Range("K54").Value = 2841.25
$ws.Range("L54").Value = 8397.8
$ws.Range("M54").Value = -2357.25
$ws.Range("N54").Value = -9365.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7030.905
$ws.Range("I31").Value = 6718.2285
$ws.Range("J31").Value = 8594.286
$ws.Range("K31").Value = 6718.2285
$ws.Range("L31").Value = 8594.286
$ws.Range("M31").Value = -6423.2285
$ws.Range("N31").Value = -9184.286

$ws.Range("H34").Value = 7030.905
$ws.Range("I34").Value = 6718.2285
$ws.Range("J34").Value = 8594.286
$ws.Range("K34").Value = 6718.2285
$ws.Range("L34").Value = 8594.286
$ws.Range("M34").Value = -6516.2285
$ws.Range("N34").Value = -8998.286

$ws.Range("H99").Value = 1787.2858
$ws.Range("I99").Value = 1268.5
$ws.Range("K99").Value = 1268.5
$ws.Range("M99").Value = 229.5

$ws.Range("H126").Value = 1787.2858
$ws.Range("I126").Value = 1268.5
$ws.Range("K126").Value = 3805.5
$ws.Range("M126").Value = -1335.5

$ws.Range("H134").Value = 3264.4285
$ws.Range("I134").Value = 3069.8667
$ws.Range("J134").Value = 3488.923
$ws.Range("K134").Value = 9209.6001
$ws.Range("L134").Value = 10466.769
$ws.Range("M134").Value = -6674.6001
$ws.Range("N134").Value = -15536.769

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 335.16666
$ws.Range("I44").Value = 166.66667
$ws.Range("J44").Value = 503.66666
$ws.Range("K44").Value = 500.00001
$ws.Range("L44").Value = 1510.99998
$ws.Range("M44").Value = -102.00001
$ws.Range("N44").Value = -2306.99998

$ws.Range("H107").Value = 918.6667
$ws.Range("I107").Value = 426.66666
$ws.Range("J107").Value = 1041.6666
$ws.Range("K107").Value = 1279.99998
$ws.Range("L107").Value = 3124.9998
$ws.Range("M107").Value = 640.00002
$ws.Range("N107").Value = -6964.9998

$ws.Range("H113").Value = 709.4915
$ws.Range("I113").Value = 710.0192
$ws.Range("J113").Value = 705.5714
$ws.Range("K113").Value = 2130.0576
$ws.Range("L113").Value = 2116.7142
$ws.Range("M113").Value = 39.94239999999991
$ws.Range("N113").Value = -6456.7142

$ws.Range("H120").Value = 7354.3184
$ws.Range("I120").Value = 20000
$ws.Range("J120").Value = 6752.143
$ws.Range("K120").Value = 60000
$ws.Range("L120").Value = 20256.429
$ws.Range("M120").Value = -55162
$ws.Range("N120").Value = -29932.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 13350
$ws.Range("J40").Value = 13350
$ws.Range("L40").Value = 13350
$ws.Range("N40").Value = -13652

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 499.25
$ws.Range("I22").Value = 475
$ws.Range("J22").Value = 507.33334
$ws.Range("K22").Value = 475
$ws.Range("L22").Value = 507.33334
$ws.Range("M22").Value = -180
$ws.Range("N22").Value = -1097.33334

$ws.Range("H27").Value = 499.25
$ws.Range("I27").Value = 475
$ws.Range("J27").Value = 507.33334
$ws.Range("K27").Value = 475
$ws.Range("L27").Value = 507.33334
$ws.Range("M27").Value = -368
$ws.Range("N27").Value = -721.33334

$ws.Range("H41").Value = 20002.75
$ws.Range("I41").Value = 19800
$ws.Range("J41").Value = 20031.715
$ws.Range("K41").Value = 19800
$ws.Range("L41").Value = 20031.715
$ws.Range("M41").Value = -19362
$ws.Range("N41").Value = -20907.715

$ws.Range("H61").Value = 3302441.8
$ws.Range("I61").Value = 3302441.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3302441.8
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3302239.8
$ws.Range("N61").ClearContents()

$ws.Range("H80").Value = 18000
$ws.Range("J80").Value = 18000
$ws.Range("L80").Value = 18000
$ws.Range("N80").Value = -20246

$ws.Range("H83").Value = 18000
$ws.Range("J83").Value = 18000
$ws.Range("L83").Value = 54000
$ws.Range("N83").Value = -65232

$ws.Range("H100").Value = 3317.2727
$ws.Range("I100").Value = 3387.7778
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 3387.7778
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -2846.7778
$ws.Range("N100").Value = -4082

$ws.Range("H102").Value = 63000
$ws.Range("J102").Value = 63000
$ws.Range("L102").Value = 63000
$ws.Range("N102").Value = -69490

$ws.Range("H113").Value = 3302441.8
$ws.Range("I113").Value = 3302441.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3302441.8
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -3300271.8
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 3256.0278
$ws.Range("I132").Value = 2860.9524
$ws.Range("J132").Value = 3809.1333
$ws.Range("K132").Value = 8582.8572
$ws.Range("L132").Value = 11427.3999
$ws.Range("M132").Value = -6052.8572
$ws.Range("N132").Value = -16487.3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2818.9443
$ws.Range("I132").Value = 2546.4375
$ws.Range("K132").Value = 7639.3125
$ws.Range("M132").Value = -5109.3125
